$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")

# --- Move the four trailing rows on Hoja1 down by 10 rows (B13:B16 -> B23:B26) ---
$v13 = $ws1.Range("B13").Value()
$v14 = $ws1.Range("B14").Value()
$v15 = $ws1.Range("B15").Value()
$v16 = $ws1.Range("B16").Value()

$ws1.Range("B13").ClearContents()
$ws1.Range("B14").ClearContents()
$ws1.Range("B15").ClearContents()
$ws1.Range("B16").ClearContents()

$ws1.Range("B23").Value = $v13
$ws1.Range("B24").Value = $v14
$ws1.Range("B25").Value = $v15
$ws1.Range("B26").Value = $v16

# --- Update Hoja1's selection ---
$null = $ws1.Range("B20").Select()

# --- Add the new "Hoja2" sheet right after "Hoja1" ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Hoja2"

# --- Populate Hoja2 with the HU (Historias de Usuario) content ---
$ws2.Range("A3").Value = "HU RELACIONADAS AL REGISTRO DE USUARIO Y ACCESO Y SOLICITUD DE COMPRA/PEDIDO"

$ws2.Range("A5").Value = "HU - 1 Yo como visitante quiero registrarme en el sitio para poder ser un usuario registrado."
$ws2.Range("A6").Value = "HU - 2  Yo como usuario registrado, quiero poder acceder al historial de pedidos."
$ws2.Range("A7").Value = "HU - 3 Yo como usuario ya registrado puede crear nuevo pedido / solicitud de compra."

$ws2.Range("A10").Value = "HU RELACIONADAS AL CRUD y A LA INTERACCION CON LA APP "

$ws2.Range("A12").Value = "HU - 4 Yo como usuario registrado puede leer y listar mis pedidos anteriores."
$ws2.Range("A13").Value = "HU - 5 Yo como usuario registrado pueda modificar mis datos personales como por ejemplo ""cambiar numero de celular"""
$ws2.Range("A14").Value = "HU - 6 Yo como usuario registrado puedo acceder a mi información personal almacenada en la DB."

$ws2.Range("A18").Value = "HU RELACIONADAS A LOS PEDIDOS Y PRODUCTOS"

$ws2.Range("A21").Value = "HU - 7  Yo como usuario registrado quiero ver la informacion del producto, disponibilidad, precio descripcion."
$ws2.Range("A22").Value = "HU - 8 Yo como usuario registrado, quiero ir agregando estos productos a un PEDIDO. "
$ws2.Range("A23").Value = "HU - 9 Yo como usuario registrado quiero revisar el DETALLE DEL PEDIDO. "
$ws2.Range("A24").Value = "HU - 10 Yo como usuario registrado, quiero ver el resumen de mi medido e ingresar la orden de compra/solicitud. "
$ws2.Range("A25").Value = "HU - 11 Yo como usuario registrado, puedo ver el resumen de mi pedido y mandarlo a imprimir y/o mandarlo por email."

# --- Column width for Hoja2 column A (closest attainable approximation of the
#     original file's bestFit-computed 109.28515625 "characters" width) ---
$ws2.Columns.Item(1).ColumnWidth = 108.5

# --- Hoja2 becomes the active/selected sheet ---
$ws2.Activate()
$null = $ws2.Range("A31").Select()
